$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: company id changes from "2" to "1" (kept as text, matching original type) ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"

# --- Row 2: updated financial metrics ---
$ws.Range("D2").Value = -0.023
$ws.Range("E2").Value = -0.103
$ws.Range("F2").Value = 0.0367
$ws.Range("G2").Value = 0.06632739842307174
$ws.Range("H2").Value = 0.06632739842307174
$ws.Range("I2").Value = 0.04842512734650906
$ws.Range("J2").Value = 0.03205389650296016
$ws.Range("K2").Value = 3323.2
$ws.Range("L2").Value = 0.02678891716572095
$ws.Range("M2").Value = 3575.1444
$ws.Range("N2").Value = 0.06266861822592745
$ws.Range("O2").Value = 1.075813793933558
$ws.Range("P2").Value = 1964.7444
$ws.Range("Q2").Value = 0.03443995624767741
$ws.Range("R2").Value = 0.5912206307173808
$ws.Range("S2").Value = 1610.4
$ws.Range("T2").Value = 0.4504433443303718
$ws.Range("U2").Value = 36271.7
$ws.Range("V2").Value = 0.6358057368830676
$ws.Range("W2").Value = 0.04376812586760511
$ws.Range("X2").Value = 0.07477979180475795
$ws.Range("Y2").Value = -0.03101166593715284
$ws.Range("Z2").Value = 1.425944211163975
$ws.Range("AA2").Value = 0.0457070681636452
$ws.Range("AB2").Value = 0.04335011110626801
$ws.Range("AC2").Value = 0.002356957057377193
$ws.Range("AD2").Value = 71707.2
$ws.Range("AF2").Value = 71707.2
$ws.Range("AG2").Value = 35435.5
$ws.Range("AH2").Value = 0.5569249026838444
$ws.Range("AI2").Value = 0.4326210583853091
$ws.Range("AJ2").Value = 0.3831531758500669
$ws.Range("AK2").Value = 0.2736777961505697
$ws.Range("AL2").Value = 1036.5
$ws.Range("AM2").Value = 1036.5
$ws.Range("AN2").Value = 9.789378839590443
$ws.Range("AO2").Value = 5.795658465991317
$ws.Range("AP2").Value = 4.837610921501707
$ws.Range("AQ2").Value = 5.795658465991317

# --- Row 3: company renamed from "April SA (ENXTPA:APR)" to "AXA SA (ENXTPA:CS)" ---
$ws.Range("B3").Value = "AXA SA (ENXTPA:CS)"

# --- Row 3: updated financial metrics (mirrors row 2, now AXA SA data) ---
$ws.Range("D3").Value = -0.023
$ws.Range("E3").Value = -0.103
$ws.Range("F3").Value = 0.0367
$ws.Range("G3").Value = 0.06632739842307174
$ws.Range("H3").Value = 0.06632739842307174
$ws.Range("I3").Value = 0.04842512734650906
$ws.Range("J3").Value = 0.03205389650296016
$ws.Range("K3").Value = 3323.2
$ws.Range("L3").Value = 0.02678891716572095
$ws.Range("M3").Value = 3575.1444
$ws.Range("N3").Value = 0.06266861822592745
$ws.Range("O3").Value = 1.075813793933558
$ws.Range("P3").Value = 1964.7444
$ws.Range("Q3").Value = 0.03443995624767741
$ws.Range("R3").Value = 0.5912206307173808
$ws.Range("S3").Value = 1610.4
$ws.Range("T3").Value = 0.4504433443303718
$ws.Range("U3").Value = 36271.7
$ws.Range("V3").Value = 0.6358057368830676
$ws.Range("W3").Value = 0.04376812586760511
$ws.Range("X3").Value = 0.07477979180475795
$ws.Range("Y3").Value = -0.03101166593715284
$ws.Range("Z3").Value = 1.425944211163975
$ws.Range("AA3").Value = 0.0457070681636452
$ws.Range("AB3").Value = 0.04335011110626801
$ws.Range("AC3").Value = 0.002356957057377193
$ws.Range("AD3").Value = 71707.2
$ws.Range("AF3").Value = 71707.2
$ws.Range("AG3").Value = 35435.5
$ws.Range("AH3").Value = 0.5569249026838444
$ws.Range("AI3").Value = 0.4326210583853091
$ws.Range("AJ3").Value = 0.3831531758500669
$ws.Range("AK3").Value = 0.2736777961505697
$ws.Range("AL3").Value = 1036.5
$ws.Range("AM3").Value = 1036.5
$ws.Range("AN3").Value = 9.789378839590443
$ws.Range("AO3").Value = 5.795658465991317
$ws.Range("AP3").Value = 4.837610921501707
$ws.Range("AQ3").Value = 5.795658465991317

# --- Remove old row 4 (AXA SA duplicate row no longer needed) ---
$ws.Rows(4).Delete()
